$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = -0.5519833787224694
$ws.Cells.Item(2, 5).Value = -31.62631796216882
$ws.Cells.Item(3, 1).Value = 0.02679286441015179
$ws.Cells.Item(3, 2).Value = 0.1175626815244795
$ws.Cells.Item(3, 3).Value = 2.427600496094155
$ws.Cells.Item(3, 4).Value = -0.5519833787224694
$ws.Cells.Item(3, 5).Value = -31.62631796216882
$ws.Cells.Item(4, 1).Value = 0.05358572882030358
$ws.Cells.Item(4, 2).Value = 0.2351253630489589
$ws.Cells.Item(4, 3).Value = 2.355200992188309
$ws.Cells.Item(4, 4).Value = -0.5505363984591253
$ws.Cells.Item(4, 5).Value = -31.54341210004048
$ws.Cells.Item(5, 1).Value = 0.08037859323045538
$ws.Cells.Item(5, 2).Value = 0.3527926821165875
$ws.Cells.Item(5, 3).Value = 2.28297167489626
$ws.Cells.Item(5, 4).Value = -0.5476401069625383
$ws.Cells.Item(5, 5).Value = -31.37746682104641
$ws.Cells.Item(6, 1).Value = 0.1071714576406072
$ws.Cells.Item(6, 2).Value = 0.4706687045231495
$ws.Cells.Item(6, 3).Value = 2.211083458931149
$ws.Cells.Item(6, 4).Value = -0.5432898272486485
$ws.Cells.Item(6, 5).Value = -31.12821415373915
$ws.Cells.Item(7, 1).Value = 0.1339643220507589
$ws.Cells.Item(7, 2).Value = 0.588856344394693
$ws.Cells.Item(7, 3).Value = 2.139708715256421
$ws.Cells.Item(7, 4).Value = -0.5374785060239089
$ws.Cells.Item(7, 5).Value = -30.79524997416677
$ws.Cells.Item(8, 1).Value = 0.1607571864609108
$ws.Cells.Item(8, 2).Value = 0.7074567678158109
$ws.Cells.Item(8, 3).Value = 2.069021999266992
$ws.Cells.Item(8, 4).Value = -0.5301966677203981
$ws.Cells.Item(8, 5).Value = -30.3780313722789
$ws.Cells.Item(9, 1).Value = 0.1875500508710625
$ws.Cells.Item(9, 2).Value = 0.8265687715341899
$ws.Cells.Item(9, 3).Value = 1.999200778830569
$ws.Cells.Item(9, 4).Value = -0.5214323522212152
$ws.Cells.Item(9, 5).Value = -29.87587308385463
$ws.Cells.Item(10, 1).Value = 0.2143429152812143
$ws.Cells.Item(10, 2).Value = 0.9462881279635648
$ws.Cells.Item(10, 3).Value = 1.930426161789286
$ws.Cells.Item(10, 4).Value = -0.5111710353989433
$ws.Cells.Item(10, 5).Value = -29.28794293769185
$ws.Cells.Item(11, 1).Value = 0.2411357796913661
$ws.Cells.Item(11, 2).Value = 1.066706887296134
$ws.Cells.Item(11, 3).Value = 1.862883622202419
$ws.Cells.Item(11, 4).Value = -0.4993955313314993
$ws.Cells.Item(11, 5).Value = -28.61325625298818
$ws.Cells.Item(12, 1).Value = 0.2679286441015179
$ws.Cells.Item(12, 2).Value = 1.187912627011252
$ws.Cells.Item(12, 3).Value = 1.796763724190174
$ws.Cells.Item(12, 4).Value = -0.4860858747730183
$ws.Cells.Item(12, 5).Value = -27.8506691054186
$ws.Cells.Item(13, 1).Value = 0.2947215085116697
$ws.Cells.Item(13, 2).Value = 1.309987638422751
$ws.Cells.Item(13, 3).Value = 1.732262841700698
$ws.Cells.Item(13, 4).Value = -0.4712191821346191
$ws.Cells.Item(13, 5).Value = -26.99887036192012
$ws.Cells.Item(14, 1).Value = 0.3215143729218215
$ws.Cells.Item(14, 2).Value = 1.433008039133984
$ws.Cells.Item(14, 3).Value = 1.669583871847275
$ws.Cells.Item(14, 4).Value = -0.4547694888616817
$ws.Cells.Item(14, 5).Value = -26.05637236309606
$ws.Cells.Item(15, 1).Value = 0.3483072373319733
$ws.Cells.Item(15, 2).Value = 1.557042799357591
$ws.Cells.Item(15, 3).Value = 1.608936938624636
$ws.Cells.Item(15, 4).Value = -0.4367075606694426
$ws.Cells.Item(15, 5).Value = -25.0215001078124
$ws.Cells.Item(16, 1).Value = 0.375100101742125
$ws.Cells.Item(16, 2).Value = 1.682152668994323
$ws.Cells.Item(16, 3).Value = 1.550540082780178
$ws.Cells.Item(16, 4).Value = -0.417000675603749
$ws.Cells.Item(16, 5).Value = -23.89237876619877
$ws.Cells.Item(17, 1).Value = 0.4018929661522769
$ws.Cells.Item(17, 2).Value = 1.808388991139053
$ws.Cells.Item(17, 3).Value = 1.494619932347522
$ws.Cells.Item(17, 4).Value = -0.3956123733121326
$ws.Cells.Item(17, 5).Value = -22.66691931393916
$ws.Cells.Item(18, 1).Value = 0.4286858305624286
$ws.Cells.Item(18, 2).Value = 1.93579238627765
$ws.Cells.Item(18, 3).Value = 1.441412346796238
$ws.Cells.Item(18, 4).Value = -0.3725021672214425
$ws.Cells.Item(18, 5).Value = -21.34280204126509
$ws.Cells.Item(19, 1).Value = 0.4554786949725804
$ws.Cells.Item(19, 2).Value = 2.064391289840305
$ws.Cells.Item(19, 3).Value = 1.391163025850361
$ws.Cells.Item(19, 4).Value = -0.3476252144966121
$ws.Cells.Item(19, 5).Value = -19.91745764298583
$ws.Cells.Item(20, 1).Value = 0.4822715593827323
$ws.Cells.Item(20, 2).Value = 2.194200323968527
$ws.Cells.Item(20, 3).Value = 1.34412807170202
$ws.Cells.Item(20, 4).Value = -0.3209319376686447
$ws.Cells.Item(20, 5).Value = -18.38804553936895
$ws.Cells.Item(21, 1).Value = 0.509064423792884
$ws.Cells.Item(21, 2).Value = 2.325218482318278
$ws.Cells.Item(21, 3).Value = 1.300574490498192
$ws.Cells.Item(21, 4).Value = -0.2923675906281415
$ws.Cells.Item(21, 5).Value = -16.75142900940111
$ws.Cells.Item(22, 1).Value = 0.5358572882030358
$ws.Cells.Item(22, 2).Value = 2.45742710444505
$ws.Cells.Item(22, 3).Value = 1.260780615487341
$ws.Cells.Item(22, 4).Value = -0.2618717602321279
$ws.Cells.Item(22, 5).Value = -15.00414663496276
$ws.Cells.Item(23, 1).Value = 0.5626501526131876
$ws.Cells.Item(23, 2).Value = 2.590787613785866
$ws.Cells.Item(23, 3).Value = 1.225036429926653
$ws.Cells.Item(23, 4).Value = -0.2293777930006249
$ws.Cells.Item(23, 5).Value = -13.14237945296124
$ws.Cells.Item(24, 1).Value = 0.5894430170233393
$ws.Cells.Item(24, 2).Value = 2.725238990460995
$ws.Cells.Item(24, 3).Value = 1.193643762578983
$ws.Cells.Item(24, 4).Value = -0.1948121342002995
$ws.Cells.Item(24, 5).Value = -11.16191308761336
$ws.Cells.Item(25, 1).Value = 0.6162358814334912
$ws.Cells.Item(25, 2).Value = 2.860694947067272
$ws.Cells.Item(25, 3).Value = 1.166916322130689
$ws.Cells.Item(25, 4).Value = -0.1580935639151607
$ws.Cells.Item(25, 5).Value = -9.058093980520432
$ws.Cells.Item(26, 1).Value = 0.643028745843643
$ws.Cells.Item(26, 2).Value = 2.997040772343073
$ws.Cells.Item(26, 3).Value = 1.145179528832243
$ws.Cells.Item(26, 4).Value = -0.1191323113441497
$ws.Cells.Item(26, 5).Value = -6.825778643658277
$ws.Cells.Item(27, 1).Value = 0.6698216102537947
$ws.Cells.Item(27, 2).Value = 3.134129804094506
$ws.Cells.Item(27, 3).Value = 1.12877009171535
$ws.Cells.Item(27, 4).Value = -0.07782902435214697
$ws.Cells.Item(27, 5).Value = -4.459274618998927
$ws.Cells.Item(28, 1).Value = 0.6966144746639465
$ws.Cells.Item(28, 2).Value = 3.27177948916313
$ws.Cells.Item(28, 3).Value = 1.118035267378297
$ws.Cells.Item(28, 4).Value = -0.03407356598126147
$ws.Cells.Item(28, 5).Value = -1.95227152368682
$ws.Cells.Item(29, 1).Value = 0.7234073390740984
$ws.Cells.Item(29, 2).Value = 3.409766984624733
$ws.Cells.Item(29, 3).Value = 1.113331720921328
$ws.Cells.Item(29, 4).Value = 0.01225639713588913
$ws.Cells.Item(29, 5).Value = 0.7022398279226775
$ws.Cells.Item(30, 1).Value = 0.7502002034842501
$ws.Cells.Item(30, 2).Value = 3.54782425106044
$ws.Cells.Item(30, 3).Value = 1.115023890339259
$ws.Cells.Item(30, 4).Value = 0.06129705818851452
$ws.Cells.Item(30, 5).Value = 3.512062730769705
$ws.Cells.Item(31, 1).Value = 0.7769930678944019
$ws.Cells.Item(31, 2).Value = 3.685632585990573
$ws.Cells.Item(31, 3).Value = 1.123481731483603
$ws.Cells.Item(31, 4).Value = 0.1132004403527653
$ws.Cells.Item(31, 5).Value = 6.485907471235867
$ws.Cells.Item(32, 1).Value = 0.8037859323045538
$ws.Cells.Item(32, 2).Value = 3.822816543959968
$ws.Cells.Item(32, 3).Value = 1.139077690219444
$ws.Cells.Item(32, 4).Value = 0.1681367432299248
$ws.Cells.Item(32, 5).Value = 9.633525768149509
$ws.Cells.Item(33, 1).Value = 0.8305787967147055
$ws.Cells.Item(33, 2).Value = 3.958937190159094
$ws.Cells.Item(33, 3).Value = 1.162182709833614
$ws.Cells.Item(33, 4).Value = 0.2262971542681229
$ws.Cells.Item(33, 5).Value = 12.96587185538434
$ws.Cells.Item(34, 1).Value = 0.8573716611248573
$ws.Cells.Item(34, 2).Value = 4.093484638154405
$ws.Cells.Item(34, 3).Value = 1.193161032754456
$ws.Cells.Item(34, 4).Value = 0.2878972393620664
$ws.Cells.Item(34, 5).Value = 16.49529674891404
$ws.Cells.Item(35, 1).Value = 0.8841645255350091
$ws.Cells.Item(35, 2).Value = 4.225869831259187
$ws.Cells.Item(35, 3).Value = 1.232363493131944
$ws.Cells.Item(35, 4).Value = 0.3531810610257711
$ws.Cells.Item(35, 5).Value = 20.23578420072905
$ws.Cells.Item(36, 1).Value = 0.9109573899451608
$ws.Cells.Item(36, 2).Value = 4.355415544314502
$ws.Cells.Item(36, 3).Value = 1.280118916719519
$ws.Cells.Item(36, 4).Value = 0.4224262189549376
$ws.Cells.Item(36, 5).Value = 24.20323950178714
$ws.Cells.Item(37, 1).Value = 0.9377502543553127
$ws.Cells.Item(37, 2).Value = 4.48134661277285
$ws.Cells.Item(37, 3).Value = 1.336723141384977
$ws.Cells.Item(37, 4).Value = 0.4959500715820613
$ws.Cells.Item(37, 5).Value = 28.41584595086318
$ws.Cells.Item(38, 1).Value = 0.9645431187654645
$ws.Cells.Item(38, 2).Value = 4.602779446056825
$ws.Cells.Item(38, 3).Value = 1.402425038265243
$ws.Cells.Item(38, 4).Value = 0.57411748604548
$ws.Cells.Item(38, 5).Value = 32.89450889506693
$ws.Cells.Item(39, 1).Value = 0.9913359831756162
$ws.Cells.Item(39, 2).Value = 4.718710963134751
$ws.Cells.Item(39, 3).Value = 1.477408740501458
$ws.Cells.Item(39, 4).Value = 0.6573505894653353
$ws.Cells.Item(39, 5).Value = 37.66341443680055
$ws.Cells.Item(40, 1).Value = 1.018128847585768
$ws.Cells.Item(40, 2).Value = 4.828007217157926
$ws.Cells.Item(40, 3).Value = 1.56177106093861
$ws.Cells.Item(40, 4).Value = 0.7461411744884274
$ws.Cells.Item(40, 5).Value = 42.75074021912123
$ws.Cells.Item(41, 1).Value = 1.04492171199592
$ws.Cells.Item(41, 2).Value = 4.929392179669064
$ws.Cells.Item(41, 3).Value = 1.655492785366055
$ws.Cells.Item(41, 4).Value = 0.841066674995241
$ws.Cells.Item(41, 5).Value = 48.1895707663286
$ws.Cells.Item(42, 1).Value = 1.071714576406072
$ws.Cells.Item(42, 2).Value = 5.021437476042833
$ws.Cells.Item(42, 3).Value = 1.75840214193048
$ws.Cells.Item(42, 4).Value = 0.9428110191720493
$ws.Cells.Item(42, 5).Value = 54.01909227698617
